# Rename the worksheet from "Sheet1" to "color"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "color"

# Update the header row text: "name" -> "DataValue", "color" -> "Color"
$ws.Range("A1").Value = "DataValue"
$ws.Range("B1").Value = "Color"

# Give the header row (A1:B1) a solid light-red/pink fill (#F4CCCC)
$ws.Range("A1:B1").Interior.Color = 13421812
